$wb = $excel.ActiveWorkbook

# Column I (2022 totals) updates for crime data added on 2022-06-19,
# applied per-sheet across the Citywide Totals, By Neighborhood summary,
# and each affected neighborhood detail sheet.

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 3111
$ws.Cells.Item(3, 9).Value = 3199
$ws.Cells.Item(4, 9).Value = 759
$ws.Cells.Item(5, 9).Value = 293
$ws.Cells.Item(6, 9).Value = 3651
$ws.Cells.Item(7, 9).Value = 11013

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(4, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 371
$ws.Cells.Item(8, 9).Value = 694
$ws.Cells.Item(10, 9).Value = 80
$ws.Cells.Item(12, 9).Value = 22
$ws.Cells.Item(13, 9).Value = 17
$ws.Cells.Item(18, 9).Value = 76
$ws.Cells.Item(20, 9).Value = 283
$ws.Cells.Item(25, 9).Value = 50
$ws.Cells.Item(27, 9).Value = 95
$ws.Cells.Item(29, 9).Value = 721
$ws.Cells.Item(31, 9).Value = 98
$ws.Cells.Item(33, 9).Value = 501
$ws.Cells.Item(36, 9).Value = 150
$ws.Cells.Item(37, 9).Value = 356
$ws.Cells.Item(42, 9).Value = 387
$ws.Cells.Item(43, 9).Value = 100
$ws.Cells.Item(47, 9).Value = 79
$ws.Cells.Item(48, 9).Value = 125
$ws.Cells.Item(49, 9).Value = 81
$ws.Cells.Item(52, 9).Value = 237
$ws.Cells.Item(54, 9).Value = 246
$ws.Cells.Item(63, 9).Value = 41
$ws.Cells.Item(65, 9).Value = 244
$ws.Cells.Item(67, 9).Value = 439
$ws.Cells.Item(68, 9).Value = 34
$ws.Cells.Item(73, 9).Value = 92
$ws.Cells.Item(76, 9).Value = 170
$ws.Cells.Item(78, 9).Value = 151
$ws.Cells.Item(79, 9).Value = 280
$ws.Cells.Item(83, 9).Value = 223
$ws.Cells.Item(85, 9).Value = 508
$ws.Cells.Item(86, 9).Value = 62
$ws.Cells.Item(88, 9).Value = 102
$ws.Cells.Item(89, 9).Value = 119
$ws.Cells.Item(90, 9).Value = 134
$ws.Cells.Item(91, 9).Value = 132
$ws.Cells.Item(94, 9).Value = 96
$ws.Cells.Item(95, 9).Value = 175
$ws.Cells.Item(98, 9).Value = 68
$ws.Cells.Item(99, 9).Value = 202
$ws.Cells.Item(101, 9).Value = 11013

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 9).Value = 130
$ws.Cells.Item(3, 9).Value = 203
$ws.Cells.Item(4, 9).Value = 31
$ws.Cells.Item(6, 9).Value = 129
$ws.Cells.Item(7, 9).Value = 508

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 63
$ws.Cells.Item(3, 9).Value = 85
$ws.Cells.Item(5, 9).Value = 8
$ws.Cells.Item(6, 9).Value = 55
$ws.Cells.Item(7, 9).Value = 237

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 9).Value = 214
$ws.Cells.Item(3, 9).Value = 191
$ws.Cells.Item(6, 9).Value = 225
$ws.Cells.Item(7, 9).Value = 694

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 9).Value = 128
$ws.Cells.Item(3, 9).Value = 107
$ws.Cells.Item(7, 9).Value = 371

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 9).Value = 25
$ws.Cells.Item(7, 9).Value = 119

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 9).Value = 112
$ws.Cells.Item(3, 9).Value = 108
$ws.Cells.Item(4, 9).Value = 26
$ws.Cells.Item(6, 9).Value = 99
$ws.Cells.Item(7, 9).Value = 356

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 9).Value = 56
$ws.Cells.Item(3, 9).Value = 69
$ws.Cells.Item(7, 9).Value = 202

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 9).Value = 105
$ws.Cells.Item(3, 9).Value = 151
$ws.Cells.Item(6, 9).Value = 152
$ws.Cells.Item(7, 9).Value = 439

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 9).Value = 31
$ws.Cells.Item(7, 9).Value = 98

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(3, 9).Value = 62
$ws.Cells.Item(4, 9).Value = 12
$ws.Cells.Item(6, 9).Value = 80
$ws.Cells.Item(7, 9).Value = 244

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 9).Value = 78
$ws.Cells.Item(3, 9).Value = 86
$ws.Cells.Item(6, 9).Value = 41
$ws.Cells.Item(7, 9).Value = 223

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 9).Value = 60
$ws.Cells.Item(3, 9).Value = 68
$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(7, 9).Value = 175

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 9).Value = 122
$ws.Cells.Item(3, 9).Value = 175
$ws.Cells.Item(7, 9).Value = 501

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(6, 9).Value = 49
$ws.Cells.Item(7, 9).Value = 81

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(6, 9).Value = 125
$ws.Cells.Item(7, 9).Value = 246

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 9).Value = 221
$ws.Cells.Item(3, 9).Value = 254
$ws.Cells.Item(5, 9).Value = 27
$ws.Cells.Item(6, 9).Value = 191
$ws.Cells.Item(7, 9).Value = 721

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(6, 9).Value = 72
$ws.Cells.Item(7, 9).Value = 125

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 9).Value = 69
$ws.Cells.Item(7, 9).Value = 170

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(3, 9).Value = 133
$ws.Cells.Item(6, 9).Value = 104
$ws.Cells.Item(7, 9).Value = 387

$ws = $wb.Worksheets.Item('Boystown')
$ws.Cells.Item(2, 9).Value = 2
$ws.Cells.Item(6, 9).Value = 17

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(3, 9).Value = 13
$ws.Cells.Item(7, 9).Value = 80

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 9).Value = 31
$ws.Cells.Item(6, 9).Value = 59
$ws.Cells.Item(7, 9).Value = 151

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(2, 9).Value = 46
$ws.Cells.Item(7, 9).Value = 132

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 9).Value = 88
$ws.Cells.Item(7, 9).Value = 280

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 9).Value = 80
$ws.Cells.Item(3, 9).Value = 82
$ws.Cells.Item(7, 9).Value = 283

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(3, 9).Value = 16
$ws.Cells.Item(7, 9).Value = 76

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 9).Value = 46
$ws.Cells.Item(7, 9).Value = 150

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 9).Value = 54
$ws.Cells.Item(7, 9).Value = 96

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(6, 9).Value = 17
$ws.Cells.Item(7, 9).Value = 50

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(2, 9).Value = 13
$ws.Cells.Item(7, 9).Value = 79

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(6, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 68

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 9).Value = 30
$ws.Cells.Item(3, 9).Value = 27
$ws.Cells.Item(7, 9).Value = 92

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 9).Value = 23
$ws.Cells.Item(3, 9).Value = 35
$ws.Cells.Item(7, 9).Value = 102

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(2, 9).Value = 23
$ws.Cells.Item(6, 9).Value = 42
$ws.Cells.Item(7, 9).Value = 95

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(4, 9).Value = 35
$ws.Cells.Item(7, 9).Value = 62

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(6, 9).Value = 49
$ws.Cells.Item(7, 9).Value = 134

$ws = $wb.Worksheets.Item('North Park')
$ws.Cells.Item(3, 9).Value = 10
$ws.Cells.Item(7, 9).Value = 34

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 9).Value = 58
$ws.Cells.Item(7, 9).Value = 100

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(3, 9).Value = 11
$ws.Cells.Item(7, 9).Value = 43

$ws = $wb.Worksheets.Item('Beverly')
$ws.Cells.Item(2, 9).Value = 4
$ws.Cells.Item(7, 9).Value = 22
